$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 3, 5 and 9 lose their "rank" (column B) value entirely - fully clear
# the cells (contents + formatting) so no stale style keeps the cell alive.
$ws.Range("B2").Clear()
$ws.Range("B3").Clear()
$ws.Range("B5").Clear()
$ws.Range("B9").Clear()

# Row 4's rank changes from 3 to 8
$ws.Range("B4").Value = 8

# Rows 12-15 gain new "rank" values (4, 3, 2, 1)
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 3
$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 1

# Selection moved to C4
$ws.Range("C4").Select()
